$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6363636363636364
$ws.Range("C2").Value = 0.7142857142857143
$ws.Range("D2").Value = 0.673076923076923
$ws.Range("E2").Value = 49

$ws.Range("B3").Value = 0.631578947368421
$ws.Range("C3").Value = 0.5454545454545454
$ws.Range("D3").Value = 0.5853658536585366
$ws.Range("E3").Value = 44

$ws.Range("B4").Value = 0.6344086021505376
$ws.Range("C4").Value = 0.6344086021505376
$ws.Range("D4").Value = 0.6344086021505376
$ws.Range("E4").Value = 0.6344086021505376

$ws.Range("B5").Value = 0.6339712918660287
$ws.Range("C5").Value = 0.6298701298701299
$ws.Range("D5").Value = 0.6292213883677298
$ws.Range("E5").Value = 93

$ws.Range("B6").Value = 0.6340999125379431
$ws.Range("C6").Value = 0.6344086021505376
$ws.Range("D6").Value = 0.6315792128144605
$ws.Range("E6").Value = 93

$ws.Range("B7").Value = 0.6296296296296297
$ws.Range("C7").Value = 0.6938775510204082
$ws.Range("D7").Value = 0.6601941747572815
$ws.Range("E7").Value = 49

$ws.Range("B8").Value = 0.6153846153846154
$ws.Range("C8").Value = 0.5454545454545454
$ws.Range("D8").Value = 0.5783132530120482
$ws.Range("E8").Value = 44

$ws.Range("B9").Value = 0.6236559139784946
$ws.Range("C9").Value = 0.6236559139784946
$ws.Range("D9").Value = 0.6236559139784946
$ws.Range("E9").Value = 0.6236559139784946

$ws.Range("B10").Value = 0.6225071225071226
$ws.Range("C10").Value = 0.6196660482374767
$ws.Range("D10").Value = 0.6192537138846648
$ws.Range("E10").Value = 93

$ws.Range("B11").Value = 0.6228900529975798
$ws.Range("C11").Value = 0.6236559139784946
$ws.Range("D11").Value = 0.6214548139315798
$ws.Range("E11").Value = 93

$ws.Range("B12").Value = 0.660377358490566
$ws.Range("C12").Value = 0.7142857142857143
$ws.Range("D12").Value = 0.6862745098039216
$ws.Range("E12").Value = 49

$ws.Range("B13").Value = 0.65
$ws.Range("C13").Value = 0.5909090909090909
$ws.Range("D13").Value = 0.6190476190476191
$ws.Range("E13").Value = 44

$ws.Range("B14").Value = 0.6559139784946236
$ws.Range("C14").Value = 0.6559139784946236
$ws.Range("D14").Value = 0.6559139784946236
$ws.Range("E14").Value = 0.6559139784946236

$ws.Range("B15").Value = 0.655188679245283
$ws.Range("C15").Value = 0.6525974025974026
$ws.Range("D15").Value = 0.6526610644257703
$ws.Range("E15").Value = 93

$ws.Range("B16").Value = 0.6554676404950295
$ws.Range("C16").Value = 0.6559139784946236
$ws.Range("D16").Value = 0.6544682389084666
$ws.Range("E16").Value = 93

$ws.Range("B17").Value = 0.64
$ws.Range("C17").Value = 0.6530612244897959
$ws.Range("D17").Value = 0.6464646464646464
$ws.Range("E17").Value = 49

$ws.Range("B18").Value = 0.6046511627906976
$ws.Range("C18").Value = 0.5909090909090909
$ws.Range("D18").Value = 0.5977011494252873
$ws.Range("E18").Value = 44

$ws.Range("B19").Value = 0.6236559139784946
$ws.Range("C19").Value = 0.6236559139784946
$ws.Range("D19").Value = 0.6236559139784946
$ws.Range("E19").Value = 0.6236559139784946

$ws.Range("B20").Value = 0.6223255813953488
$ws.Range("C20").Value = 0.6219851576994434
$ws.Range("D20").Value = 0.6220828979449669
$ws.Range("E20").Value = 93

$ws.Range("B21").Value = 0.6232758189547386
$ws.Range("C21").Value = 0.6236559139784946
$ws.Range("D21").Value = 0.6233937446395733
$ws.Range("E21").Value = 93

$ws.Range("B22").Value = 0.6491228070175439
$ws.Range("C22").Value = 0.7551020408163265
$ws.Range("D22").Value = 0.6981132075471698
$ws.Range("E22").Value = 49

$ws.Range("B23").Value = 0.6666666666666666
$ws.Range("C23").Value = 0.5454545454545454
$ws.Range("D23").Value = 0.6
$ws.Range("E23").Value = 44

$ws.Range("B24").Value = 0.6559139784946236
$ws.Range("C24").Value = 0.6559139784946236
$ws.Range("D24").Value = 0.6559139784946236
$ws.Range("E24").Value = 0.6559139784946236

$ws.Range("B25").Value = 0.6578947368421053
$ws.Range("C25").Value = 0.650278293135436
$ws.Range("D25").Value = 0.6490566037735849
$ws.Range("E25").Value = 93

$ws.Range("B26").Value = 0.6574231277117526
$ws.Range("C26").Value = 0.6559139784946236
$ws.Range("D26").Value = 0.651694055589369
$ws.Range("E26").Value = 93
